$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B so the old numeric metric columns
# (previously B:E) shift right to C:F, turning the previous column A
# (plain numbers) into column B, and freeing up column A for the new
# "Metodo" (method name) label column.
$ws.Columns("B:B").Insert()

# ---- Text labels, written in the same order the author typed them so the
# shared-strings table comes out in the same sequence as the target file:
# Metodo, SMARTER, Fuzzy, TOPSIS, GRA, CODAS, MABAC, VIKOR, PROMETHEE II,
# Rx, Ry, CL, Entropia, SSIM ----
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# ---- Data rows (values unchanged, shifted one column right into B:F) ----
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.79222740479055287
$ws.Range("E2").Value = 7.9035305827081093
$ws.Range("F2").Value = 0.61578413807230903

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0.35479685648929304
$ws.Range("E3").Value = 7.753686463944506
$ws.Range("F3").Value = 0.68234173833005929

$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0.97880382474159178
$ws.Range("E4").Value = 7.8877697597406691
$ws.Range("F4").Value = 0.64031770130013688

$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0.79222740479055287
$ws.Range("E5").Value = 7.9035305827081093
$ws.Range("F5").Value = 0.61578413807230903

$ws.Range("B6").Value = 95
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 7.9084888412453891
$ws.Range("F6").Value = 0.58485659856240635

$ws.Range("B7").Value = 95
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 7.9084888412453891
$ws.Range("F7").Value = 0.58485659856240635

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0.36898499644653843
$ws.Range("E8").Value = 7.6236212627024926
$ws.Range("F8").Value = 0.71355915992509844

$ws.Range("B9").Value = 95
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 7.9084888412453891
$ws.Range("F9").Value = 0.58485659856240635

# Widen the new "Metodo" column so the longer method names (e.g.
# "PROMETHEE II") are fully visible, matching the column now being
# split out of the old A:B bestFit group. Columns B:C and D:F are left
# untouched so they keep inheriting their original best-fit widths.
$ws.Columns("A:A").ColumnWidth = 12.416666876524687

Write-Output "edit applied"
